# Update "Forecast Comparison" sheet with a new Week_Start_Date column
# and corrected forecast output (per commit: "Update with Correct Forecast output").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new blank column before column B ("ASIN"). This shifts the
# existing columns B:I (ASIN..is_holiday_week) to C:J, preserving all
# of their existing values automatically.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week_Start_Date values for each forecast week row (2-17).
$weekStartDates = @(
    "2024-12-08",
    "2024-12-15",
    "2024-12-22",
    "2024-12-29",
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23"
)

# Updated Week labels (dropping the leading zero: "W01" -> "W1", etc.)
$weekLabels = @(
    "W1", "W2", "W3", "W4", "W5", "W6", "W7", "W8", "W9",
    "W10", "W11", "W12", "W13", "W14", "W15", "W16"
)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weekLabels[$i]

    # Force the date column to be stored as plain text (not converted to
    # a date serial number) so it matches "YYYY-MM-DD" string values.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $weekStartDates[$i]
}
